# Daily attendance processing - 2025-12-24 10:33:37
#
# This script reproduces, via Excel COM automation, the results of the
# daily attendance-processing run against the "Session Analysis Results"
# sheet:
#   1. The "Recorded By" text "dnasr281@gmail.com, System" is re-ordered
#      to "System, dnasr281@gmail.com" everywhere it appears.
#   2. Six General-Surgery B1 sessions dated 24/12/2025 (today, per the
#      commit) rolled from "Pending" to "Not Recorded" because they were
#      never recorded once their date passed - status text + row
#      highlight (yellow -> pink/red) change, and the Status column is
#      widened to fit the longer label.
#   3. The per-group Missing/Pending counters and the sheet-wide
#      Missing/Pending Sessions totals are recalculated to reflect the
#      six newly-missed sessions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Recorded-By text re-ordered: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$recordedByCells = "G2,G3,G20,G21,G22,G24,G39,G40,G41,G43,G58,G59,G60,G62,G77,G78,G95,G96,G113,G114,G131,G132,G149,G150,G167,G168,G169,G171,G186,G187,G188,G190,G205,G206,G207,G209"
$ws.Range($recordedByCells).Value = "System, dnasr281@gmail.com"

# --- 2) Sessions that passed today (24/12/2025) without being recorded: Pending -> Not Recorded
$missedRows = @(25, 44, 63, 172, 191, 210)
foreach ($rowNum in $missedRows) {
    $rowRange = $ws.Range("A" + $rowNum + ":I" + $rowNum)

    # Status label
    $ws.Range("I" + $rowNum).Value = "Not Recorded"

    # Row highlight: yellow (Pending) -> pink (Not Recorded); font/alignment unchanged
    $rowRange.Interior.Color = 12695295
    $rowRange.Font.Color = 0
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108
}

# --- 3) Status column is now wider to fit "Not Recorded"
$ws.Columns.Item(9).ColumnWidth = 13.166666666666666

# --- 4) Sheet-wide Missing/Pending Sessions totals
$ws.Range("L7").Value = 6
$ws.Range("L8").Value = 144

# --- 5) Per-group Missing/Pending counters for the six groups with a session
#         that rolled from Pending to Missing today
$groupRows = @(16, 17, 18, 24, 25, 26)
foreach ($rowNum in $groupRows) {
    $ws.Range("P" + $rowNum).Value = 1
    $ws.Range("Q" + $rowNum).Value = 12
}
